# Update the "Obrigatorio" column (E) for rows 2-12 from "N" to "S"
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

for ($row = 2; $row -le 12; $row++) {
    $ws.Cells.Item($row, 5).Value = "S"
}
